$wb = $excel.ActiveWorkbook

# Sheet "Central America" -> header A1 should read "Central America"
$wsCA = $wb.Worksheets.Item("Central America")
$wsCA.Range("A1").Value = "Central America"

# Sheet "Eastern Europe" -> header A1 should read "Eastern Europe"
$wsEE = $wb.Worksheets.Item("Eastern Europe")
$wsEE.Range("A1").Value = "Eastern Europe"

# Select cell A2 on each sheet, then activate Eastern Europe last so it becomes the active tab
$wsCA.Range("A2").Select()
$wsEE.Activate()
$wsEE.Range("A2").Select()
